# Add support for tabs to 3d model page, support axis customization

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("files")

# --- Add the new "x axis"/"y axis"/"z axis" header columns on the files sheet ---
$ws1.Range("F1").Value = "x axis"
$ws1.Range("G1").Value = "y axis"
$ws1.Range("H1").Value = "z axis"
$ws1.Range("F1:H1").Font.Bold = $true

# --- Fill in the default axis mapping (x/y/z) for every data row ---
for ($r = 2; $r -le 10; $r++) {
    $ws1.Cells.Item($r, 6).Value = "x"
    $ws1.Cells.Item($r, 7).Value = "y"
    $ws1.Cells.Item($r, 8).Value = "z"
}

# --- Add the new "notes" worksheet (placed after "files") describing the
#     coordinate-axis mapping ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "notes"

$ws2.Columns.Item(1).ColumnWidth = 63.6640625

$ws2.Range("A1").Value = "Plotly Dash Coordinate Axes"
$ws2.Range("A1").Font.Bold = $true

$notesText = "The Plotly Dash 3d surface capability displays data using a right-handed coordinate system with the z axis extending in the vertical direction. `nIf you want to display your data differently, you can map the axes in your data to Plotly's axes using the vector columns in the files tab. This will cause the data to be labelled according to Plotly's axes labelling. So if you map your z axis to Plotly's y axis, it will be labelled as the y axis in the resulting plot.`nValid mapping values are x, y, and z."
$ws2.Range("A2").Value = $notesText
$ws2.Range("A2").VerticalAlignment = -4108
$ws2.Range("A2").WrapText = $true
$ws2.Rows.Item(2).RowHeight = 153

[void]$ws2.Range("A5").Select()

# --- Restore "files" as the active sheet/selection ---
$ws1.Activate()
[void]$ws1.Range("C26").Select()
